# Fix typo in the decaf LLVM slide: "llvm::getGlobalContext()" -> "TheContext"
#
# Slide 8 has a code listing textbox (Shape 3) whose paragraph reads:
#     llvm::getGlobalContext(), 
# split across four runs:
#     [err=1] "llvm"  [plain] "::"  [err=1] "getGlobalContext"  [plain] "(), "
# It should become two runs:
#     [err=1] "TheContext"  [plain] ", "

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(8)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

$oldCall = "llvm::getGlobalContext(), "

$text = $tr.Text
$start = $text.IndexOf($oldCall)
if ($start -lt 0) {
    throw "Could not find target text '$oldCall' in slide 8 textbox"
}

# First run: "llvm" -> "TheContext" (keeps its own rPr, incl. err="1")
$runLlvm = $tr.Characters($start + 1, 4)
$runLlvm.Text = "TheContext"

# Remainder: "::getGlobalContext(), " -> ", " (keeps the "::" run's rPr, no err)
$text = $tr.Text
$remainder = "::getGlobalContext(), "
$start2 = $text.IndexOf($remainder)
if ($start2 -lt 0) {
    throw "Could not find remainder text '$remainder' in slide 8 textbox"
}
$runRest = $tr.Characters($start2 + 1, $remainder.Length)
$runRest.Text = ", "
